$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: volume number 22 -> 23 (A8, shared rich-text string "Volume 30   Number  22")
$ws.Range("A8").Value = "Volume 30   Number  23"

# Header: week covering dates (C9, shared rich-text string)
$ws.Range("C9").Value = "Report Covering the Week  6/5/2023  Through  6/11/2023"

# Weekly crime statistics table (rows 14-30)
$ws.Range("D14").Value = 7
$ws.Range("E14").Value = -28.571428571428
$ws.Range("F14").Value = 26
$ws.Range("G14").Value = 37
$ws.Range("H14").Value = -29.729729729729
$ws.Range("I14").Value = 168
$ws.Range("J14").Value = 192
$ws.Range("K14").Value = -12.5
$ws.Range("L14").Value = -18.840579710144
$ws.Range("M14").Value = -17.241379310344
$ws.Range("N14").Value = -79.856115107913
$ws.Range("C15").Value = 22
$ws.Range("D15").Value = 38
$ws.Range("E15").Value = -42.105263157894
$ws.Range("F15").Value = 108
$ws.Range("G15").Value = 147
$ws.Range("H15").Value = -26.530612244898
$ws.Range("I15").Value = 650
$ws.Range("J15").Value = 718
$ws.Range("K15").Value = -9.470752089136
$ws.Range("L15").Value = 5.008077544426
$ws.Range("M15").Value = 17.753623188405
$ws.Range("N15").Value = -54.672245467224
$ws.Range("C16").Value = 323
$ws.Range("D16").Value = 372
$ws.Range("E16").Value = -13.172043010752
$ws.Range("F16").Value = 1241
$ws.Range("G16").Value = 1424
$ws.Range("H16").Value = -12.851123595505
$ws.Range("I16").Value = 6861
$ws.Range("J16").Value = 7172
$ws.Range("K16").Value = -4.336307863915
$ws.Range("L16").Value = 34.899724734565
$ws.Range("M16").Value = -13.096896770107
$ws.Range("N16").Value = -81.252561685384
$ws.Range("C17").Value = 542
$ws.Range("D17").Value = 614
$ws.Range("E17").Value = -11.726384364820
$ws.Range("F17").Value = 2184
$ws.Range("G17").Value = 2285
$ws.Range("H17").Value = -4.420131291028
$ws.Range("I17").Value = 11671
$ws.Range("J17").Value = 10994
$ws.Range("K17").Value = 6.157904311442
$ws.Range("L17").Value = 28.182317408017
$ws.Range("M17").Value = 60.404068169323
$ws.Range("N17").Value = -32.345950959364
$ws.Range("C18").Value = 214
$ws.Range("D18").Value = 295
$ws.Range("E18").Value = -27.457627118644
$ws.Range("F18").Value = 981
$ws.Range("G18").Value = 1154
$ws.Range("H18").Value = -14.991334488734
$ws.Range("I18").Value = 6196
$ws.Range("J18").Value = 6797
$ws.Range("K18").Value = -8.842136236574
$ws.Range("L18").Value = 21.681068342498
$ws.Range("M18").Value = -20.113460546673
$ws.Range("N18").Value = -85.702418312719
$ws.Range("C19").Value = 1043
$ws.Range("D19").Value = 1029
$ws.Range("E19").Value = 1.360544217687
$ws.Range("F19").Value = 3942
$ws.Range("G19").Value = 3897
$ws.Range("H19").Value = 1.154734411085
$ws.Range("I19").Value = 21549
$ws.Range("J19").Value = 21727
$ws.Range("K19").Value = -0.819257145487
$ws.Range("L19").Value = 50.114942528735
$ws.Range("M19").Value = 38.436335603237
$ws.Range("N19").Value = -40.020040637960
$ws.Range("C20").Value = 352
$ws.Range("D20").Value = 271
$ws.Range("E20").Value = 29.889298892988
$ws.Range("F20").Value = 1214
$ws.Range("G20").Value = 978
$ws.Range("H20").Value = 24.130879345603
$ws.Range("I20").Value = 6645
$ws.Range("J20").Value = 5650
$ws.Range("K20").Value = 17.610619469026
$ws.Range("L20").Value = 76.775738228252
$ws.Range("M20").Value = 51.125767568796
$ws.Range("N20").Value = -86.618470337105
$ws.Range("C21").Value = 2501
$ws.Range("D21").Value = 2626
$ws.Range("E21").Value = -4.760091393754
$ws.Range("F21").Value = 9696
$ws.Range("G21").Value = 9922
$ws.Range("H21").Value = -2.277766579318
$ws.Range("I21").Value = 53740
$ws.Range("J21").Value = 53250
$ws.Range("K21").Value = 0.920187793427
$ws.Range("L21").Value = 40.595976244669
$ws.Range("M21").Value = 23.129797227632
$ws.Range("N21").Value = -70.957159919367
$ws.Range("C22").Value = 41
$ws.Range("D22").Value = 44
$ws.Range("E22").Value = -6.818181818181
$ws.Range("F22").Value = 168
$ws.Range("G22").Value = 188
$ws.Range("H22").Value = -10.638297872340
$ws.Range("I22").Value = 958
$ws.Range("J22").Value = 1040
$ws.Range("K22").Value = -7.884615384615
$ws.Range("L22").Value = 44.277108433734
$ws.Range("M22").Value = 3.121636167922
$ws.Range("C23").Value = 117
$ws.Range("D23").Value = 163
$ws.Range("E23").Value = -28.220858895705
$ws.Range("F23").Value = 469
$ws.Range("G23").Value = 523
$ws.Range("H23").Value = -10.325047801147
$ws.Range("I23").Value = 2692
$ws.Range("J23").Value = 2580
$ws.Range("K23").Value = 4.341085271317
$ws.Range("L23").Value = 17.811816192560
$ws.Range("M23").Value = 57.981220657277
$ws.Range("C24").Value = 2189
$ws.Range("D24").Value = 2347
$ws.Range("E24").Value = -6.731998295696
$ws.Range("F24").Value = 8872
$ws.Range("G24").Value = 9236
$ws.Range("H24").Value = -3.941100043308
$ws.Range("I24").Value = 47628
$ws.Range("J24").Value = 48440
$ws.Range("K24").Value = -1.676300578034
$ws.Range("L24").Value = 40.338263892981
$ws.Range("M24").Value = 40.173053152039
$ws.Range("C25").Value = 869
$ws.Range("D25").Value = 917
$ws.Range("E25").Value = -5.234460196292
$ws.Range("F25").Value = 3651
$ws.Range("G25").Value = 3554
$ws.Range("H25").Value = 2.729319077096
$ws.Range("I25").Value = 18959
$ws.Range("J25").Value = 17995
$ws.Range("K25").Value = 5.357043623228
$ws.Range("L25").Value = 33.608174770965
$ws.Range("M25").Value = -5.620270808442
$ws.Range("C26").Value = 38
$ws.Range("D26").Value = 65
$ws.Range("E26").Value = -41.538461538461
$ws.Range("F26").Value = 203
$ws.Range("G26").Value = 246
$ws.Range("H26").Value = -17.479674796748
$ws.Range("I26").Value = 1086
$ws.Range("J26").Value = 1188
$ws.Range("K26").Value = -8.585858585858
$ws.Range("L26").Value = 4.624277456647
$ws.Range("C27").Value = 106
$ws.Range("D27").Value = 119
$ws.Range("E27").Value = -10.924369747899
$ws.Range("F27").Value = 454
$ws.Range("G27").Value = 477
$ws.Range("H27").Value = -4.821802935010
$ws.Range("I27").Value = 2307
$ws.Range("J27").Value = 2219
$ws.Range("K27").Value = 3.965750337990
$ws.Range("L27").Value = 17.524197656648
$ws.Range("C28").Value = 17
$ws.Range("D28").Value = 30
$ws.Range("E28").Value = -43.333333333333
$ws.Range("F28").Value = 86
$ws.Range("G28").Value = 130
$ws.Range("H28").Value = -33.846153846153
$ws.Range("I28").Value = 483
$ws.Range("J28").Value = 649
$ws.Range("K28").Value = -25.57781201849
$ws.Range("L28").Value = -31.779661016949
$ws.Range("M28").Value = -29.489051094890
$ws.Range("N28").Value = -80.500605571255
$ws.Range("C29").Value = 15
$ws.Range("D29").Value = 27
$ws.Range("E29").Value = -44.444444444444
$ws.Range("F29").Value = 77
$ws.Range("G29").Value = 113
$ws.Range("H29").Value = -31.858407079646
$ws.Range("I29").Value = 412
$ws.Range("J29").Value = 554
$ws.Range("K29").Value = -25.631768953068
$ws.Range("L29").Value = -33.868378812199
$ws.Range("M29").Value = -26.296958855098
$ws.Range("N29").Value = -81.639928698752
$ws.Range("D30").Value = 19
$ws.Range("E30").Value = -84.210526315789
$ws.Range("F30").Value = 35
$ws.Range("G30").Value = 52
$ws.Range("H30").Value = -32.692307692307
$ws.Range("I30").Value = 213
$ws.Range("J30").Value = 310
$ws.Range("K30").Value = -31.290322580645
$ws.Range("L30").Value = -16.796875

